$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B26 (ID_Padre) from "5.1." to "5.1.1."
$ws.Range("B26").Value = "5.1.1."

# Update E28 (Tipo) from empty to "Carpeta"
$ws.Range("E28").Value = "Carpeta"

# Update B31 (ID_Padre) from "5.2" to "5.2."
$ws.Range("B31").Value = "5.2."

# Remove the highlight style from A30/A31 (revert to normal "center, no fill" style)
$ws.Range("A29").Copy()
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("A29").Copy()
$ws.Range("A31").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Delete the now-empty trailing row 32 (also shrinks the Excel table range)
$ws.Rows("32").Delete()

# Move the active selection to B32 (first empty row after the shifted table)
$ws.Range("B32").Select()
